$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "20.582.72"
$ws.Range("E2").Value2 = "  +0.27%  "

$ws.Range("D3").Value2 = "1.480.31"
$ws.Range("E3").Value2 = "  +0.65%  "

$ws.Range("E4").Value2 = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.9717"
$ws.Range("E5").Value2 = "  +2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "280.16"
$ws.Range("E6").Value2 = "  -0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.3660"
$ws.Range("E7").Value2 = "  -1.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3082"
$ws.Range("E8").Value2 = "  -3.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "40.14"
$ws.Range("E9").Value2 = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "1.061"
$ws.Range("E10").Value2 = "  +0.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.06667"
$ws.Range("E11").Value2 = "  -0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "1.007"
$ws.Range("E12").Value2 = "  +0.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.524"
$ws.Range("E13").Value2 = "  -1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "18.11"
$ws.Range("E14").Value2 = "  -0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "6.212"
$ws.Range("E15").Value2 = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.9716"
$ws.Range("E16").Value2 = "  +3.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.00001029"
$ws.Range("E17").Value2 = "  -0.42%  "

$ws.Range("D18").Value2 = "1.483.76"
$ws.Range("E18").Value2 = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.05955"
$ws.Range("E19").Value2 = "  +3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "69.85"
$ws.Range("E20").Value2 = "  -3.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "5.502"
$ws.Range("E21").Value2 = "  -3.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "14.52"
$ws.Range("E22").Value2 = "  -1.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "11.07"
$ws.Range("E23").Value2 = "  -1.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "2.272"
$ws.Range("E24").Value2 = "  +0.00%  "

$ws.Range("D25").Value2 = "20.634.25"
$ws.Range("E25").Value2 = "  -0.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "142.16"
$ws.Range("E26").Value2 = "  +2.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.120"
$ws.Range("E27").Value2 = "  -8.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "17.29"
$ws.Range("E28").Value2 = "  -1.71%  "

$ws.Range("D29").Value2 = "1.645.97"
$ws.Range("E29").Value2 = "  +0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "114.11"
$ws.Range("E30").Value2 = "  +0.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "3.979"
$ws.Range("E31").Value2 = "  +0.67%  "

$ws.Range("B32").Value2 = "Filecoin"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "5.041"
$ws.Range("E32").Value2 = "  -5.13%  "

$ws.Range("B33").Value2 = "ImmutableX"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.8223"
$ws.Range("E33").Value2 = "  -3.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.08009"
$ws.Range("E34").Value2 = "  +2.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.545"
$ws.Range("E35").Value2 = "  -4.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.224"
$ws.Range("E36").Value2 = "  +9.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.05777"
$ws.Range("E37").Value2 = "  -4.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.756"
$ws.Range("E38").Value2 = "  -3.65%  "

$ws.Range("B39").Value2 = "Frax"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.9711"
$ws.Range("E39").Value2 = "  +1.43%  "

$ws.Range("B40").Value2 = "VeChain"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.02047"
$ws.Range("E40").Value2 = "  -1.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "7.646"
$ws.Range("E41").Value2 = "  +1.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "10.43"
$ws.Range("E42").Value2 = "  -2.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1882"
$ws.Range("E43").Value2 = "  -1.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.5313"
$ws.Range("E44").Value2 = "  -1.84%  "

$ws.Range("B45").Value2 = "EnergySwap"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "12.33"
$ws.Range("E45").Value2 = "  -1.08%  "

$ws.Range("B46").Value2 = "PancakeSwap"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.543"
$ws.Range("E46").Value2 = "  -1.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "118.84"
$ws.Range("E47").Value2 = "  -2.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.5206"
$ws.Range("E48").Value2 = "  -2.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.825"
$ws.Range("E49").Value2 = "  -0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.06504"
$ws.Range("E50").Value2 = "  +0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.9899"
$ws.Range("E51").Value2 = "  -0.22%  "
